# Basic Git Commands.docx - apply commit "Re #1 review question, create starting query"

$d = $word.ActiveDocument

function Set-CalibriFont($range) {
    $range.Font.Name = "Calibri"
    $range.Font.NameAscii = "Calibri"
    $range.Font.NameBi = "Calibri"
    $range.Font.Size = 16
    $range.Font.SizeBi = 16
}

# ---------------------------------------------------------------------------
# 1. Fill in the previously-empty second paragraph (right after the title)
#    with the new intro sentence, then add four brand-new paragraphs after
#    it (still above the empty paragraph that precedes the table).
# ---------------------------------------------------------------------------

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.Text = "The commands are executed in a command window. VS, VS Code and GitHub desktop all have the ability to open command windows."
Set-CalibriFont $r2

# --- paragraph: "To execute you must have Git installed. If you have Github Desktop installed then you already have Git installed." ---
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$r3.Text = "To execute you must have Git installed. If you have Github Desktop installed then you already have Git installed."
Set-CalibriFont $r3

# --- paragraph: "Typical syntax: git command options" ---
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.Text = "Typical syntax: git command options"
Set-CalibriFont $r4

# --- paragraph: "To trace and prepare your repository changes for committing you can use: git add . or git add -a" ---
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5.Text = "To trace and prepare your repository changes for committing you can use: git add . or git add -a"
Set-CalibriFont $r5

# --- paragraph: "For help on your commands you can use: git command -help ... or git command -help ..." ---
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
$r6.Text = "For help on your commands you can use: git command -help (which places a summary of options on the command window) or git command " + [char]0x2013 + "help (which opens your help in a browser window)"
Set-CalibriFont $r6

Write-Host "Intro paragraphs done. Paragraph count: " $d.Paragraphs.Count
